$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 9.427210000000001
$ws.Range("H2").Value = 28.28163
$ws.Range("I2").Value = 0.2188083857550241
$ws.Range("J2").Value = 0.2188083857550241
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.055614
$ws.Range("N2").Value = 0.166842
$ws.Range("O2").Value = 0.003173711121411028
$ws.Range("P2").Value = 0.003173711121411028
$ws.Range("Q2").Value = 0.52428485694
$ws.Range("R2").Value = 4.71856371246
$ws.Range("S2").Value = 0.0006944346073287145
$ws.Range("T2").Value = 0.0006944346073287144
# Row 3
$ws.Range("G3").Value = 9.427210000000001
$ws.Range("H3").Value = 28.28163
$ws.Range("I3").Value = 0.2188083857550241
$ws.Range("J3").Value = 0.2188083857550241
$ws.Range("O3").Value = 0.01439065952479895
$ws.Range("P3").Value = 0.01439065952479895
$ws.Range("Q3").Value = 2.377281542523333
$ws.Range("R3").Value = 21.39553388271
$ws.Range("S3").Value = 0.003148796980571422
$ws.Range("T3").Value = 0.003148796980571422
# Row 4
$ws.Range("G4").Value = 9.427210000000001
$ws.Range("H4").Value = 28.28163
$ws.Range("I4").Value = 0.2188083857550241
$ws.Range("J4").Value = 0.2188083857550241
$ws.Range("M4").Value = 4.376294333333333
$ws.Range("N4").Value = 13.128883
$ws.Range("O4").Value = 0.2497409644382361
$ws.Range("P4").Value = 0.2497409644382361
$ws.Range("Q4").Value = 41.25624570214333
$ws.Range("R4").Value = 371.30621131929
$ws.Range("S4").Value = 0.05464541728563332
$ws.Range("T4").Value = 0.05464541728563332
# Row 5
$ws.Range("G5").Value = 9.427210000000001
$ws.Range("H5").Value = 28.28163
$ws.Range("I5").Value = 0.2188083857550241
$ws.Range("J5").Value = 0.2188083857550241
$ws.Range("M5").Value = 12.83925333333333
$ws.Range("N5").Value = 38.51776
$ws.Range("O5").Value = 0.732694664915554
$ws.Range("P5").Value = 0.7326946649155539
$ws.Range("Q5").Value = 121.0383374165333
$ws.Range("R5").Value = 1089.3450367488
$ws.Range("S5").Value = 0.1603197368814907
$ws.Range("T5").Value = 0.1603197368814906
# Row 6
$ws.Range("I6").Value = 0.3808887290954196
$ws.Range("J6").Value = 0.3808887290954196
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.055614
$ws.Range("N6").Value = 0.166842
$ws.Range("O6").Value = 0.003173711121411028
$ws.Range("P6").Value = 0.003173711121411028
$ws.Range("Q6").Value = 0.9126441482339999
$ws.Range("R6").Value = 8.213797334105999
$ws.Range("S6").Value = 0.001208830795550246
$ws.Range("T6").Value = 0.001208830795550245
# Row 7
$ws.Range("I7").Value = 0.3808887290954196
$ws.Range("J7").Value = 0.3808887290954196
$ws.Range("O7").Value = 0.01439065952479895
$ws.Range("P7").Value = 0.01439065952479895
$ws.Range("S7").Value = 0.005481240017245569
$ws.Range("T7").Value = 0.005481240017245568
# Row 8
$ws.Range("I8").Value = 0.3808887290954196
$ws.Range("J8").Value = 0.3808887290954196
$ws.Range("M8").Value = 4.376294333333333
$ws.Range("N8").Value = 13.128883
$ws.Range("O8").Value = 0.2497409644382361
$ws.Range("P8").Value = 0.2497409644382361
$ws.Range("Q8").Value = 71.81643856342433
$ws.Range("R8").Value = 646.347947070819
$ws.Range("S8").Value = 0.09512351854794412
$ws.Range("T8").Value = 0.09512351854794411
# Row 9
$ws.Range("I9").Value = 0.3808887290954196
$ws.Range("J9").Value = 0.3808887290954196
$ws.Range("M9").Value = 12.83925333333333
$ws.Range("N9").Value = 38.51776
$ws.Range("O9").Value = 0.732694664915554
$ws.Range("P9").Value = 0.7326946649155539
$ws.Range("Q9").Value = 210.6963969928533
$ws.Range("R9").Value = 1896.26757293568
$ws.Range("S9").Value = 0.2790751397346797
$ws.Range("T9").Value = 0.2790751397346796
# Row 10
$ws.Range("G10").Value = 7.213061
$ws.Range("H10").Value = 21.639183
$ws.Range("I10").Value = 0.1674173200514808
$ws.Range("J10").Value = 0.1674173200514808
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.055614
$ws.Range("N10").Value = 0.166842
$ws.Range("O10").Value = 0.003173711121411028
$ws.Range("P10").Value = 0.003173711121411028
$ws.Range("Q10").Value = 0.4011471744539999
$ws.Range("R10").Value = 3.610324570086
$ws.Range("S10").Value = 0.0005313342105642141
$ws.Range("T10").Value = 0.000531334210564214
# Row 11
$ws.Range("G11").Value = 7.213061
$ws.Range("H11").Value = 21.639183
$ws.Range("I11").Value = 0.1674173200514808
$ws.Range("J11").Value = 0.1674173200514808
$ws.Range("O11").Value = 0.01439065952479895
$ws.Range("P11").Value = 0.01439065952479895
$ws.Range("Q11").Value = 1.818934422845667
$ws.Range("R11").Value = 16.370409805611
$ws.Range("S11").Value = 0.002409245651415157
$ws.Range("T11").Value = 0.002409245651415157
# Row 12
$ws.Range("G12").Value = 7.213061
$ws.Range("H12").Value = 21.639183
$ws.Range("I12").Value = 0.1674173200514808
$ws.Range("J12").Value = 0.1674173200514808
$ws.Range("M12").Value = 4.376294333333333
$ws.Range("N12").Value = 13.128883
$ws.Range("O12").Value = 0.2497409644382361
$ws.Range("P12").Value = 0.2497409644382361
$ws.Range("Q12").Value = 31.56647798028767
$ws.Range("R12").Value = 284.098301822589
$ws.Range("S12").Value = 0.04181096297332165
$ws.Range("T12").Value = 0.04181096297332165
# Row 13
$ws.Range("G13").Value = 7.213061
$ws.Range("H13").Value = 21.639183
$ws.Range("I13").Value = 0.1674173200514808
$ws.Range("J13").Value = 0.1674173200514808
$ws.Range("M13").Value = 12.83925333333333
$ws.Range("N13").Value = 38.51776
$ws.Range("O13").Value = 0.732694664915554
$ws.Range("P13").Value = 0.7326946649155539
$ws.Range("Q13").Value = 92.61031748778666
$ws.Range("R13").Value = 833.4928573900801
$ws.Range("S13").Value = 0.1226657772161798
$ws.Range("T13").Value = 0.1226657772161797
# Row 14
$ws.Range("G14").Value = 10.03371566666667
$ws.Range("H14").Value = 30.101147
$ws.Range("I14").Value = 0.2328855650980756
$ws.Range("J14").Value = 0.2328855650980756
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.055614
$ws.Range("N14").Value = 0.166842
$ws.Range("O14").Value = 0.003173711121411028
$ws.Range("P14").Value = 0.003173711121411028
$ws.Range("Q14").Value = 0.558015063086
$ws.Range("R14").Value = 5.022135567774001
$ws.Range("S14").Value = 0.0007391115079678545
$ws.Range("T14").Value = 0.0007391115079678545
# Row 15
$ws.Range("G15").Value = 10.03371566666667
$ws.Range("H15").Value = 30.101147
$ws.Range("I15").Value = 0.2328855650980756
$ws.Range("J15").Value = 0.2328855650980756
$ws.Range("O15").Value = 0.01439065952479895
$ws.Range("P15").Value = 0.01439065952479895
$ws.Range("Q15").Value = 2.530225491666556
$ws.Range("R15").Value = 22.772029424999
$ws.Range("S15").Value = 0.003351376875566808
$ws.Range("T15").Value = 0.003351376875566808
# Row 16
$ws.Range("G16").Value = 10.03371566666667
$ws.Range("H16").Value = 30.101147
$ws.Range("I16").Value = 0.2328855650980756
$ws.Range("J16").Value = 0.2328855650980756
$ws.Range("M16").Value = 4.376294333333333
$ws.Range("N16").Value = 13.128883
$ws.Range("O16").Value = 0.2497409644382361
$ws.Range("P16").Value = 0.2497409644382361
$ws.Range("Q16").Value = 43.91049301431122
$ws.Range("R16").Value = 395.194437128801
$ws.Range("S16").Value = 0.05816106563133702
$ws.Range("T16").Value = 0.05816106563133701
# Row 17
$ws.Range("G17").Value = 10.03371566666667
$ws.Range("H17").Value = 30.101147
$ws.Range("I17").Value = 0.2328855650980756
$ws.Range("J17").Value = 0.2328855650980756
$ws.Range("M17").Value = 12.83925333333333
$ws.Range("N17").Value = 38.51776
$ws.Range("O17").Value = 0.732694664915554
$ws.Range("P17").Value = 0.7326946649155539
$ws.Range("Q17").Value = 128.8254173189689
$ws.Range("R17").Value = 1159.42875587072
$ws.Range("S17").Value = 0.1706340110832039
$ws.Range("T17").Value = 0.1706340110832039
